$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.494.50'
$ws.Range('D3').Value = '1.601.97'
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '215.17'
$ws.Range('E5').Value = '  +2.38%  '
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  +1.88%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '24.03'
$ws.Range('E8').Value = '  +9.31%  '
$ws.Range('D9').Value = '0.252'
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('D10').Value = '0.0602'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('D12').Value = '1.832.82'
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('D13').Value = '1.604.21'
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('D14').Value = '3.80'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('E15').Value = '  +3.73%  '
$ws.Range('D16').Value = '28.534.65'
$ws.Range('E16').Value = '  +5.18%  '
$ws.Range('D17').Value = '63.41'
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('D18').Value = '232.72'
$ws.Range('E18').Value = '  +7.64%  '
$ws.Range('D19').Value = '7.56'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').Value = '0.0₃0711'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('E23').Value = '  +2.76%  '
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').Value = '152.66'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  +2.48%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('D34').Value = '1.425.44'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('E36').Value = '  -4.22%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = '0.545'
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('D40').Value = '2.52'
$ws.Range('E40').Value = '  +8.26%  '
$ws.Range('D41').Value = '0.824'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('E42').Value = '  -2.76%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  +7.06%  '
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('D46').Value = '65.00'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('D47').Value = '1.742.82'
$ws.Range('E47').Value = '  +2.81%  '
$ws.Range('D48').Value = '87.63'
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  +7.97%  '
$ws.Range('E51').Value = '  +0.70%  '
